$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark after
#    "...rall amount of testing performed"
# -------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# -------------------------------------------------------------------
# 2. "Code-driven testing" paragraph: "that are returned" -> "returned"
# -------------------------------------------------------------------
$d.Content.Find.Execute(
    "validate that the results that are returned are correct.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "validate that the results returned are correct.", 2)

# -------------------------------------------------------------------
# 3. "The intent of this paper..." paragraph: insert "multiple" before
#    "GUI testing framework tools" and re-create the "_GoBack" bookmark
#    right after the newly inserted word (matching the target run
#    layout: "...focus on " | "multiple" | <bookmark> | " GUI testing...")
# -------------------------------------------------------------------

# Mark the boundary right before the word "the" that will be replaced.
$rLeft = $d.Content
$rLeft.Find.Execute("The intent of this paper is to focus on ", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$rLeft.Collapse(0)
$rLeft.Bookmarks.Add("boundaryA")

# Mark the boundary right after the word "the" (i.e. right before " GUI").
$rRight = $d.Content
$rRight.Find.Execute("The intent of this paper is to focus on the", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
$rRight.Collapse(0)
$rRight.Bookmarks.Add("boundaryB")

# Replace the word sitting exactly between the two boundaries ("the") with "multiple".
$bmA = $d.Bookmarks.Item("boundaryA")
$bmB = $d.Bookmarks.Item("boundaryB")
$rMid = $d.Range($bmA.End, $bmB.Start)
$rMid.Text = "multiple"

# Drop the temporary left boundary - it was only needed to keep "multiple" in its own run.
$d.Bookmarks.Item("boundaryA").Delete()

# Turn the temporary right boundary into the real "_GoBack" bookmark.
$bmB2 = $d.Bookmarks.Item("boundaryB")
$posStart = $bmB2.Start
$posEnd = $bmB2.End
$bmB2.Delete()
$d.Range($posStart, $posEnd).Bookmarks.Add("_GoBack")
